$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix swapped latitude/longitude values in row 4 (J4 and K4 were transposed)
$jVal = $ws.Range("J4").Value2
$kVal = $ws.Range("K4").Value2
$ws.Range("J4").Value2 = $kVal
$ws.Range("K4").Value2 = $jVal

# Update the active selection to K4, matching the saved view state
$ws.Range("K4").Select()
